$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.191837543572632
$ws.Range("C2").Value = 0.03121349681488539
$ws.Range("D2").Value = 0.01447854751909361
$ws.Range("F2").Value = 4.308721678539314
$ws.Range("G2").Value = 0.002645403414384401
$ws.Range("J2").Value = 0.217323224979566
$ws.Range("K2").Value = 1.594393377898427
$ws.Range("L2").Value = 0.3583678508754957
$ws.Range("N2").Value = 4.231800862614222

$ws.Range("B3").Value = 2.154203712022706
$ws.Range("C3").Value = 0.02719316049150677
$ws.Range("D3").Value = 0.01440869971301595
$ws.Range("F3").Value = 4.296154326153413
$ws.Range("G3").Value = 0.002649876345347199
$ws.Range("J3").Value = 0.2177141506405746
$ws.Range("K3").Value = 1.55572932451571
$ws.Range("L3").Value = 0.3555524453458148
$ws.Range("N3").Value = 4.239546622931925

$ws.Range("B4").Value = 2.132284954798621
$ws.Range("C4").Value = 0.02472898349374475
$ws.Range("D4").Value = 0.01437448850994016
$ws.Range("F4").Value = 4.289990300900115
$ws.Range("G4").Value = 0.002652769077497049
$ws.Range("J4").Value = 0.2180138728566128
$ws.Range("K4").Value = 1.532916606725195
$ws.Range("L4").Value = 0.353995462863459
$ws.Range("N4").Value = 4.245147570029204

$ws.Range("B5").Value = 2.12365189884602
$ws.Range("C5").Value = 0.0237258355697918
$ws.Range("D5").Value = 0.01436274393960524
$ws.Range("F5").Value = 4.287868759566933
$ws.Range("G5").Value = 0.002653984805001135
$ws.Range("J5").Value = 0.2181510420446102
$ws.Range("K5").Value = 1.523853346751423
$ws.Range("L5").Value = 0.3534042061199614
$ws.Range("N5").Value = 4.247642375203313

$ws.Range("B6").Value = 2.122236450689996
$ws.Range("C6").Value = 0.02355932363396107
$ws.Range("D6").Value = 0.01436092689200485
$ws.Range("F6").Value = 4.287540055312149
$ws.Range("G6").Value = 0.002654188908882891
$ws.Range("J6").Value = 0.218174727244282
$ws.Range("K6").Value = 1.522362477667258
$ws.Range("L6").Value = 0.3533086407620303
$ws.Range("N6").Value = 4.248069459321073

$ws.Range("B7").Value = 2.132167315497298
$ws.Range("C7").Value = 0.02471545062726932
$ws.Range("D7").Value = 0.01437432120345861
$ws.Range("F7").Value = 4.289960108571236
$ws.Range("G7").Value = 0.002652785323446039
$ws.Range("J7").Value = 0.2180156618861631
$ws.Range("K7").Value = 1.532793432632076
$ws.Range("L7").Value = 0.3539873138735885
$ws.Range("N7").Value = 4.245180356123569

$ws.Range("B8").Value = 2.178614775203215
$ws.Range("C8").Value = 0.0298263475440308
$ws.Range("D8").Value = 0.01445267237014036
$ws.Range("F8").Value = 4.304066238716416
$ws.Range("G8").Value = 0.002646915380943282
$ws.Range("J8").Value = 0.2174456386405339
$ws.Range("K8").Value = 1.580869538105389
$ws.Range("L8").Value = 0.3573614954943451
$ws.Range("N8").Value = 4.234296145813033

$ws.Range("B9").Value = 2.279130779985735
$ws.Range("C9").Value = 0.0398864773145533
$ws.Range("D9").Value = 0.01467452253948309
$ws.Range("F9").Value = 4.34405206172147
$ws.Range("G9").Value = 0.002636560115697793
$ws.Range("J9").Value = 0.216800691013308
$ws.Range("K9").Value = 1.682512461426711
$ws.Range("L9").Value = 0.3653389708807566
$ws.Range("N9").Value = 4.219663458800625

$ws.Range("B10").Value = 2.358744925782332
$ws.Range("C10").Value = 0.04730596797961084
$ws.Range("D10").Value = 0.01487829429667542
$ws.Range("F10").Value = 4.3809597685771
$ws.Range("G10").Value = 0.002629649034796738
$ws.Range("J10").Value = 0.216614249406426
$ws.Range("K10").Value = 1.761703519417182
$ws.Range("L10").Value = 0.3720286807632363
$ws.Range("N10").Value = 4.213014822421258

$ws.Range("B11").Value = 2.396219344950282
$ws.Range("C11").Value = 0.05068863214711428
$ws.Range("D11").Value = 0.01497968633751867
$ws.Range("F11").Value = 4.399389556146616
$ws.Range("G11").Value = 0.002626654718478335
$ws.Range("J11").Value = 0.2165916703437105
$ws.Range("K11").Value = 1.798715648199277
$ws.Range("L11").Value = 0.3752518540436114
$ws.Range("N11").Value = 4.210883406657572

$ws.Range("B12").Value = 2.410590873473609
$ws.Range("C12").Value = 0.05197071959496213
$ws.Range("D12").Value = 0.01501931713021065
$ws.Range("F12").Value = 4.406604528254121
$ws.Range("G12").Value = 0.002625542233867077
$ws.Range("J12").Value = 0.2165920536413068
$ws.Range("K12").Value = 1.812873481381274
$ws.Range("L12").Value = 0.3764982394978205
$ws.Range("N12").Value = 4.210204909534923

$ws.Range("B13").Value = 2.407487668209626
$ws.Range("C13").Value = 0.05169454707882437
$ws.Range("D13").Value = 0.01501072719531749
$ws.Range("F13").Value = 4.405040155897296
$ws.Range("G13").Value = 0.00262578087751454
$ws.Range("J13").Value = 0.216591574018473
$ws.Range("K13").Value = 1.809818014971995
$ws.Range("L13").Value = 0.3762286595309519
$ws.Range("N13").Value = 4.210345312214315

$ws.Range("B14").Value = 2.397398076439174
$ws.Range("C14").Value = 0.05079408689914544
$ws.Range("D14").Value = 0.01498292210224861
$ws.Range("F14").Value = 4.399978405384047
$ws.Range("G14").Value = 0.002626562765655088
$ws.Range("J14").Value = 0.2165915229257251
$ws.Range("K14").Value = 1.799877572227672
$ws.Range("L14").Value = 0.3753538772995313
$ws.Range("N14").Value = 4.210825007319286

$ws.Range("B15").Value = 2.391241451516748
$ws.Range("C15").Value = 0.05024268028283529
$ws.Range("D15").Value = 0.01496605119170269
$ws.Range("F15").Value = 4.396908677037231
$ws.Range("G15").Value = 0.002627044478042739
$ws.Range("J15").Value = 0.2165926545747112
$ws.Range("K15").Value = 1.793807278707277
$ws.Range("L15").Value = 0.3748214117667317
$ws.Range("N15").Value = 4.211135591039948

$ws.Range("B16").Value = 2.356321184514059
$ws.Range("C16").Value = 0.04708505959857234
$ws.Range("D16").Value = 0.01487184181176104
$ws.Range("F16").Value = 4.379788357146509
$ws.Range("G16").Value = 0.002629847720614701
$ws.Range("J16").Value = 0.2166169759630563
$ws.Range("K16").Value = 1.759304574429081
$ws.Range("L16").Value = 0.3718216569994155
$ws.Range("N16").Value = 4.213172103501137

$ws.Range("B17").Value = 2.335220777408949
$ws.Range("C17").Value = 0.04514993043071058
$ws.Range("D17").Value = 0.01481626406216563
$ws.Range("F17").Value = 4.36970580887629
$ws.Range("G17").Value = 0.002631605648930597
$ws.Range("J17").Value = 0.2166478247411803
$ws.Range("K17").Value = 1.73839133721151
$ws.Range("L17").Value = 0.3700274737305165
$ws.Range("N17").Value = 4.214650318147804

$ws.Range("B18").Value = 2.32320275003184
$ws.Range("C18").Value = 0.04403759725087752
$ws.Range("D18").Value = 0.01478511650448766
$ws.Range("F18").Value = 4.364060976370013
$ws.Range("G18").Value = 0.002632630847202182
$ws.Range("J18").Value = 0.2166714268372729
$ws.Range("K18").Value = 1.726455557968393
$ws.Range("L18").Value = 0.3690124510280981
$ws.Range("N18").Value = 4.215584597027529

$ws.Range("B19").Value = 2.319153984206764
$ws.Range("C19").Value = 0.04366109886332481
$ws.Range("D19").Value = 0.01477471163752142
$ws.Range("F19").Value = 4.362176244297373
$ws.Range("G19").Value = 0.002632980384519361
$ws.Range("J19").Value = 0.2166804248207157
$ws.Range("K19").Value = 1.722430271234202
$ws.Range("L19").Value = 0.3686716932622858
$ws.Range("N19").Value = 4.215915357267534

$ws.Range("B20").Value = 2.337454701498928
$ws.Range("C20").Value = 0.04535585509894702
$ws.Range("D20").Value = 0.01482209571949866
$ws.Range("F20").Value = 4.370763134405053
$ws.Range("G20").Value = 0.002631417057576592
$ws.Range("J20").Value = 0.2166439345991407
$ws.Range("K20").Value = 1.740607965076833
$ws.Range("L20").Value = 0.3702167142713364
$ws.Range("N20").Value = 4.214484259550446

$ws.Range("B21").Value = 2.400356728199881
$ws.Range("C21").Value = 0.05105854224100881
$ws.Range("D21").Value = 0.01499105570665193
$ws.Range("F21").Value = 4.401458757755108
$ws.Range("G21").Value = 0.002626332526120714
$ws.Range("J21").Value = 0.2165912956040259
$ws.Range("K21").Value = 1.802793462568218
$ws.Range("L21").Value = 0.37561012109947
$ws.Range("N21").Value = 4.210680616809455

$ws.Range("B22").Value = 2.442520541171234
$ws.Range("C22").Value = 0.05479228902717637
$ws.Range("D22").Value = 0.01510867756303824
$ws.Range("F22").Value = 4.422895850068215
$ws.Range("G22").Value = 0.002623134162488488
$ws.Range("J22").Value = 0.2166089571810161
$ws.Range("K22").Value = 1.844263857610741
$ws.Range("L22").Value = 0.3792856187514246
$ws.Range("N22").Value = 4.208944530401254

$ws.Range("B23").Value = 2.419920522041025
$ws.Range("C23").Value = 0.0527988816314604
$ws.Range("D23").Value = 0.01504524681419994
$ws.Range("F23").Value = 4.411328542353999
$ws.Range("G23").Value = 0.002624829818119801
$ws.Range("J23").Value = 0.2165947720933232
$ws.Range("K23").Value = 1.82205447352564
$ws.Range("L23").Value = 0.377310170875802
$ws.Range("N23").Value = 4.209802437129511

$ws.Range("B24").Value = 2.33644439281511
$ws.Range("C24").Value = 0.04526275594395202
$ws.Range("D24").Value = 0.01481945671926255
$ws.Range("F24").Value = 4.370284644824324
$ws.Range("G24").Value = 0.00263150227437929
$ws.Range("J24").Value = 0.2166456750575705
$ws.Range("K24").Value = 1.739605555017789
$ws.Range("L24").Value = 0.370131107296416
$ws.Range("N24").Value = 4.214559071640522

$ws.Range("B25").Value = 2.250927230313835
$ws.Range("C25").Value = 0.03716036159285352
$ws.Range("D25").Value = 0.01460728799186128
$ws.Range("F25").Value = 4.331913749322837
$ws.Range("G25").Value = 0.002639238554392964
$ws.Range("J25").Value = 0.2169246339595112
$ws.Range("K25").Value = 1.654224388930118
$ws.Range("L25").Value = 0.3630352729392428
$ws.Range("N25").Value = 4.210345312214315

Write-Output "applied 380kV case updates"
